$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) D20: "dose" -> "intervention dose"
$ws.Range("D20").Value = "intervention dose"

# 2) Insert a new row at row 32 (pushes existing row 32.. down by one)
$ws.Rows.Item(32).Insert()

# Copy the formatting (fill style) used by the "Proposed" GMHO rows (style index 4)
# from what is now row 41 (the "mental health intervention content" row) onto the
# newly inserted, still-blank row 32.
$ws.Range("A41:V41").Copy()
$ws.Range("A32:V32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 3) Populate the new row 32 with the "intervention dose" term
$ws.Range("A32").Value = "GMHO:0000242"
$ws.Range("B32").Value = "intervention dose"
$ws.Range("C32").Value = "An intervention attribute that is its amount."
$ws.Range("D32").Value = "intervention attribute"
$ws.Range("H32").Value = "drug dose and number of intervention techniques (e.g., action plans, social support)"
$ws.Range("I32").Value = "Intervention dose can capture the intensity of any intervention, including social, psychological and pharmacological interventions.`nIntervention dose refers to the intensity or amount of an intervention, which is about the content of an intervention. In contrast, an intervention's schedule of delivery is about the temporal aspects of the intervention. There is some overlap between the entities ""intervention dose"" and ""intervention schedule of delivery"". For example, a more frequent intervention schedule with specific doses would suggest a higher overall dose of the intervention. However, ""intervention dose"" can capture concepts that are not directly relevant to the temporal organisation of the intervention, such as overall number of intervention strategies and drug dose.`n"
$ws.Range("S32").Value = "Proposed"

# Re-fit the row height (setting the long I32 text nudges an automatic
# "customHeight" onto the row; AutoFit recalculates and drops the override
# since wrap text isn't enabled on this sheet, matching the other rows).
$ws.Rows.Item(32).AutoFit()
